$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 3: "Mini Waterpomp x1" purchase (first test pump order) ---
$ws.Range("A3").Value = "Mini Waterpomp x1"
$ws.Range("C3").Value = 4.37

# --- New "Datum" column (D) ---
$ws.Range("D1").Value = "Datum"

# D3 gets a real date value (built-in short-date format -> numFmtId 14)
$d = Get-Date -Year 2025 -Month 3 -Day 10 -Hour 0 -Minute 0 -Second 0
$ws.Range("D3").Value = $d
$ws.Range("D3").NumberFormat = "mm-dd-yy"

# D2 gets the date typed as text, with a long-date display format
$ws.Range("D2").Value = "19/09/2025"
$ws.Range("D2").NumberFormat = "[`$-F800]dddd, mmmm dd, yyyy"

# --- Move "Totaal" column from E to F ---
$ws.Range("F1").Value = $ws.Range("E1").Value2
$ws.Range("F2").NumberFormat = $ws.Range("E2").NumberFormat
$ws.Range("F2").Formula = $ws.Range("E2").Formula
$ws.Range("E1").Clear()
$ws.Range("E2").Clear()

# --- Column D width ---
$ws.Columns.Item(4).ColumnWidth = 10.08984375

# --- Selection, matching the saved state ---
$ws.Range("P19").Select()
